$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indexes used: B=2 C=3 D=4 E=5 F=6 I=9 J=10 K=11 L=12 M=13 N=14 (G=7 stays at 1, H=8 has no data)
$colIndexes = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @{
    2 = @(1.02,1.008075152743659,1.028374011839279,1.010687419312362,1.022662078082575,1.030497932126107,1.013343089779908,1.031191015190591,1.0135568571881,1.025495799515954,1.008593558734579)
    3 = @(1.02,1.009066242985069,1.028823185320063,1.011528180645183,1.023932323243426,1.03059556605843,1.013965721747002,1.031449079130273,1.014201645205447,1.026571493132391,1.008802624249196)
    4 = @(1.02,1.009707708500172,1.029111076978838,1.012072743032263,1.024750946688625,1.030655357900006,1.014368234242181,1.031612588806549,1.014618762380534,1.027263695956164,1.008937680483503)
    5 = @(1.02,1.009977419908773,1.029231444902873,1.012301804239871,1.025094303997144,1.030679682762889,1.014537361210864,1.031680494244084,1.014794093258957,1.027553777807771,1.008994404577092)
    6 = @(1.02,1.010022707948034,1.029251616328428,1.012340272065855,1.025151908768942,1.03068371937468,1.014565753130685,1.031691846894798,1.014823530598245,1.027602429856731,1.009003925658117)
    7 = @(1.02,1.009711312242725,1.029112687944387,1.012075803260358,1.024755537757599,1.030655686121202,1.014370494476182,1.031613499442317,1.014621105259736,1.027267575660579,1.008938438644204)
    8 = @(1.02,1.008410062454691,1.02852638079061,1.010971447726031,1.023092049068479,1.030531627714169,1.013553588069444,1.031278947365684,1.013774787610294,1.025860130797662,1.008664259529617)
    9 = @(1.02,1.006118344307147,1.027472250161201,1.009029542544614,1.020135391487081,1.030287169581864,1.012111254022426,1.030662899831217,1.012282686675195,1.02335060527296,1.008179417263183)
    10 = @(1.02,1.004591362254832,1.026755527626665,1.007737731780969,1.018147174587711,1.030106899902852,1.011147796748991,1.030234484062072,1.011287443360672,1.021657795829191,1.007855049779483)
    11 = @(1.02,1.003930355491933,1.026441893307249,1.007179031082851,1.017282181509308,1.030024754799628,1.010730159018035,1.030044795824873,1.010856373672895,1.020920090980698,1.007714325819567)
    12 = @(1.02,1.003684855744885,1.026324903551513,1.006971604542183,1.016960269692132,1.029993629606143,1.010574961470965,1.029973710656089,1.010696236861647,1.020645366552564,1.007662014048544)
    13 = @(1.02,1.003737515032252,1.026350020482722,1.007016093719059,1.017029348720555,1.030000333779533,1.010608254962992,1.029988986990607,1.010730587579397,1.020704327927282,1.007673236938809)
    14 = @(1.02,1.003910061845842,1.026432232922478,1.007161883103043,1.017255584730899,1.030022194476502,1.010717331733902,1.030038932673754,1.010843137087593,1.020897396621752,1.007710002540797)
    15 = @(1.02,1.004016377276115,1.026482821585771,1.007251721985019,1.017394894677476,1.03003558240062,1.010784528480953,1.03006962288813,1.010912480127514,1.021016258773905,1.007732649659373)
    16 = @(1.02,1.004635234990551,1.026776273393788,1.007774824880294,1.018204495191791,1.030112265595209,1.011175504414058,1.03024698509493,1.011316049459196,1.02170665562929,1.007864383465707)
    17 = @(1.02,1.005023477728558,1.026959468972069,1.008103130952846,1.018711242233537,1.030159273558201,1.01142063161988,1.030357121465034,1.011569165052211,1.022138462555479,1.007946944160334)
    18 = @(1.02,1.005249951177237,1.027066006402557,1.008294690171396,1.019006425463499,1.030186298079939,1.011563566334414,1.030420958594083,1.011716791243675,1.022389874112876,1.007995074295579)
    19 = @(1.02,1.005327175841419,1.027102279003287,1.008360017732753,1.019107008577743,1.030195445804343,1.011612295943544,1.030442656903313,1.011767125946243,1.022475521953818,1.008011480995928)
    20 = @(1.02,1.004981821095209,1.026939846640301,1.008067900212379,1.018656913809377,1.030154270837246,1.011394336324931,1.030345346603613,1.011542009369465,1.022092180735281,1.007938088880608)
    21 = @(1.02,1.003859250351586,1.026408036973385,1.007118949006285,1.017188980851429,1.030015773954213,1.010685213229255,1.030024242196463,1.01080999456801,1.020840562254509,1.007699177108544)
    22 = @(1.02,1.00315360501026,1.026070821521077,1.006522883139697,1.016262474809301,1.029925150295622,1.010238964212224,1.029818727095951,1.0103496417668,1.020049523299537,1.007548728759573)
    23 = @(1.02,1.003527665990984,1.026249854864848,1.006838814093312,1.016753971038632,1.029973527264826,1.010475566733607,1.029928017506868,1.010593693516067,1.020469256632238,1.007628506541719)
    24 = @(1.02,1.005000643878735,1.026948714108358,1.008083819270126,1.018681463701221,1.030156532570217,1.011406218172111,1.030350668403599,1.011554279890064,1.022113094899857,1.007942090281002)
    25 = @(1.02,1.00671066125412,1.027747240151959,1.00953108142368,1.020902770449074,1.030353422593306,1.012484467622265,1.030825293903405,1.01266852141279,1.02400286545531,1.008304961848062)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $colIndexes.Length; $i++) {
        $col = $colIndexes[$i]
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
